$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1571.110000049326
$ws.Range("C2").Value = 274.5400000030443
$ws.Range("D2").Value = 2788.830000029995
$ws.Range("B3").Value = 2383.240000027599
$ws.Range("C3").Value = 297.3199999984293
$ws.Range("D3").Value = 3318.210000016494
$ws.Range("B4").Value = 3167.01000002484
$ws.Range("C4").Value = 292.0400000030752
$ws.Range("D4").Value = 3856.950000006812
$ws.Range("B5").Value = 3929.380000025288
$ws.Range("C5").Value = 292.4500000063812
$ws.Range("D5").Value = 4345.590000004549
$ws.Range("B6").Value = 4505.530000018184
$ws.Range("C6").Value = 311.4300000009518
$ws.Range("D6").Value = 4718.390000000028
$ws.Range("B7").Value = 4997.76000001644
$ws.Range("C7").Value = 312.0700000017562
$ws.Range("D7").Value = 5030.869999999372
$ws.Range("B8").Value = 5410.030000016259
$ws.Range("C8").Value = 331.6100000020058
$ws.Range("D8").Value = 5275.739999995401
$ws.Range("B9").Value = 5845.5200000103
$ws.Range("C9").Value = 343.7600000025757
$ws.Range("D9").Value = 5553.899999997995
$ws.Range("B10").Value = 6268.930000013829
$ws.Range("C10").Value = 347.7900000005106
$ws.Range("D10").Value = 5778.159999993276
$ws.Range("B11").Value = 6685.410000010121
$ws.Range("C11").Value = 354.2300000000048
$ws.Range("D11").Value = 6036.049999996207
$ws.Range("B12").Value = 6968.70000000934
$ws.Range("C12").Value = 357.7699999983905
$ws.Range("D12").Value = 6182.399999994025
$ws.Range("B13").Value = 7315.880000013184
$ws.Range("C13").Value = 373.4699999995082
$ws.Range("D13").Value = 6337.779999992428
$ws.Range("B14").Value = 7642.340000007943
$ws.Range("C14").Value = 385.2800000009273
$ws.Range("D14").Value = 6505.66999999364
$ws.Range("B15").Value = 7942.840000003797
$ws.Range("C15").Value = 395.6700000015693
$ws.Range("D15").Value = 6623.369999998818
$ws.Range("B16").Value = 8211.29000000339
$ws.Range("C16").Value = 401.2499999987305
$ws.Range("D16").Value = 6731.089999994469
$ws.Range("B17").Value = 8421.110000003007
$ws.Range("C17").Value = 412.3299999969155
$ws.Range("D17").Value = 6804.669999988783
$ws.Range("B18").Value = 8687.100000003398
$ws.Range("C18").Value = 413.3100000000913
$ws.Range("D18").Value = 6891.139999989607
$ws.Range("B19").Value = 8971.030000005254
$ws.Range("C19").Value = 439.259999998176
$ws.Range("D19").Value = 6975.359999986497
$ws.Range("B20").Value = 9256.210000000277
$ws.Range("C20").Value = 457.7999999929814
$ws.Range("D20").Value = 7071.049999983691
$ws.Range("B21").Value = 9430.740000004758
$ws.Range("C21").Value = 450.6400000024637
$ws.Range("D21").Value = 7093.659999957395
$ws.Range("B22").Value = 9565.119999998285
$ws.Range("C22").Value = 454.6500000005836
$ws.Range("D22").Value = 7081.369999981678
$ws.Range("B23").Value = 9660.709999993989
$ws.Range("C23").Value = 461.4300000003938
$ws.Range("D23").Value = 7030.679999993657
$ws.Range("B24").Value = 9942.999999999667
$ws.Range("C24").Value = 466.5200000015402
$ws.Range("D24").Value = 7097.599999974625
$ws.Range("B25").Value = 10234.00000000109
$ws.Range("C25").Value = 489.9599999979274
$ws.Range("D25").Value = 7165.579999947106
$ws.Range("B26").Value = 10403.49000000138
$ws.Range("C26").Value = 504.0899999975995
$ws.Range("D26").Value = 7137.80999996567
$ws.Range("B27").Value = 10621.34999999932
$ws.Range("C27").Value = 515.7499999977515
$ws.Range("D27").Value = 7128.009999975996
$ws.Range("B28").Value = 10734.14000000031
$ws.Range("C28").Value = 523.4199999977933
$ws.Range("D28").Value = 7071.06999997977
$ws.Range("B29").Value = 10908.66999999769
$ws.Range("C29").Value = 538.2199999988134
$ws.Range("D29").Value = 7044.279999981088
$ws.Range("B30").Value = 11092.19000000006
$ws.Range("C30").Value = 543.8399999975765
$ws.Range("D30").Value = 7003.649999986912
$ws.Range("B31").Value = 11159.10999999905
$ws.Range("C31").Value = 544.749999994508
$ws.Range("D31").Value = 6882.629999988557
$ws.Range("B32").Value = 11235.19000000273
$ws.Range("C32").Value = 555.5199999975032
$ws.Range("D32").Value = 6781.669999990558
$ws.Range("B33").Value = 11479.58999999977
$ws.Range("C33").Value = 568.4399999981766
$ws.Range("D33").Value = 6753.159999991508
$ws.Range("B34").Value = 11613.32999999987
$ws.Range("C34").Value = 586.6099999963463
$ws.Range("D34").Value = 6670.659999990923
$ws.Range("B35").Value = 11764.54000000327
$ws.Range("C35").Value = 605.3699999955072
$ws.Range("D35").Value = 6583.619999989382
$ws.Range("B36").Value = 11969.87999999765
$ws.Range("C36").Value = 611.9199999967831
$ws.Range("D36").Value = 6517.909999991484
$ws.Range("B37").Value = 12049.89999999633
$ws.Range("C37").Value = 613.1499999974591
$ws.Range("D37").Value = 6396.819999992181
$ws.Range("B38").Value = 12130.57999999643
$ws.Range("C38").Value = 624.2700000025008
$ws.Range("D38").Value = 6270.549999993863
$ws.Range("B39").Value = 12359.99000000052
$ws.Range("C39").Value = 651.5099999980665
$ws.Range("D39").Value = 6186.699999991491
$ws.Range("B40").Value = 12588.01999999735
$ws.Range("C40").Value = 669.1499999948815
$ws.Range("D40").Value = 6114.609999994512
$ws.Range("B41").Value = 12602.5599999957
$ws.Range("C41").Value = 672.2899999989572
$ws.Range("D41").Value = 5944.189999993511
$ws.Range("B42").Value = 12707.61999999675
$ws.Range("C42").Value = 674.1499999963195
$ws.Range("D42").Value = 5807.61999999308
$ws.Range("B43").Value = 12739.23999999427
$ws.Range("C43").Value = 687.8499999930369
$ws.Range("D43").Value = 5640.059999997088
$ws.Range("B44").Value = 12832.78999999597
$ws.Range("C44").Value = 689.2899999943718
$ws.Range("D44").Value = 5494.319999995717
$ws.Range("B45").Value = 13008.69999999173
$ws.Range("C45").Value = 709.4299999928879
$ws.Range("D45").Value = 5369.499999994586
$ws.Range("B46").Value = 13115.37999998869
$ws.Range("C46").Value = 721.1799999967396
$ws.Range("D46").Value = 5200.52999999441
$ws.Range("B47").Value = 13150.98999999522
$ws.Range("C47").Value = 724.7999999955713
$ws.Range("D47").Value = 5031.069999993471
$ws.Range("B48").Value = 13395.45999998197
$ws.Range("C48").Value = 743.3099999950385
$ws.Range("D48").Value = 4921.509999999733
$ws.Range("B49").Value = 13370.26999999225
$ws.Range("C49").Value = 748.7999999896318
$ws.Range("D49").Value = 4720.519999999476
$ws.Range("B50").Value = 13333.14000000043
$ws.Range("C50").Value = 743.880000002035
$ws.Range("D50").Value = 4534.58000000141
$ws.Range("B51").Value = 13445.13999999474
$ws.Range("C51").Value = 755.8299999987887
$ws.Range("D51").Value = 4373.049999997904
$ws.Range("B52").Value = 13551.29999999237
$ws.Range("C52").Value = 767.4399999974082
$ws.Range("D52").Value = 4199.790000001665
$ws.Range("B53").Value = 13608.52999998706
$ws.Range("C53").Value = 771.2599999963095
$ws.Range("D53").Value = 4001.119999999608
$ws.Range("B54").Value = 13726.71999998878
$ws.Range("C54").Value = 783.4699999957804
$ws.Range("D54").Value = 3840.499999999764
$ws.Range("B55").Value = 13918.67999997712
$ws.Range("C55").Value = 794.05999999166
$ws.Range("D55").Value = 3678.789999997231
$ws.Range("B56").Value = 13742.87999999287
$ws.Range("C56").Value = 781.8199999985792
$ws.Range("D56").Value = 3454.109999999766
$ws.Range("B57").Value = 13944.9999999824
$ws.Range("C57").Value = 802.4799999931125
$ws.Range("D57").Value = 3275.180000007429
$ws.Range("B58").Value = 13988.45999997337
$ws.Range("C58").Value = 805.7099999845256
$ws.Range("D58").Value = 3083.880000006923
$ws.Range("B59").Value = 13970.2299999893
$ws.Range("C59").Value = 800.9499999959012
$ws.Range("D59").Value = 2890.300000010958
$ws.Range("B60").Value = 14025.54999997833
$ws.Range("C60").Value = 814.3699999850752
$ws.Range("D60").Value = 2695.140000009903
$ws.Range("B61").Value = 14078.76999997624
$ws.Range("C61").Value = 807.7899999927205
$ws.Range("D61").Value = 2487.400000016846
$ws.Range("B62").Value = 14218.23999996147
$ws.Range("C62").Value = 820.3499999898225
$ws.Range("D62").Value = 2303.57000002623
